$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.852.52"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.221.42"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.81"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.30"
$ws.Range("E6").Value = "  +5.13%  "
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.86"
$ws.Range("E10").Value = "  +6.41%  "
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.20"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.563.44"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.09"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.216.52"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.731"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.807.69"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0884"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("E21").Value = "  +7.37%  "
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.63"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.45"
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.77"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.96"
$ws.Range("E31").Value = "  +4.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.88"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("E35").Value = "  +4.11%  "
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.80"
$ws.Range("E37").Value = "  +7.38%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.98"
$ws.Range("E39").Value = "  +5.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0998"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("E42").Value = "  +5.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.063.33"
$ws.Range("E43").Value = "  +8.82%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0268"
$ws.Range("E44").Value = "  +3.91%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.96"
$ws.Range("E45").Value = "  +10.85%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.83"
$ws.Range("E46").Value = "  +11.04%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.434.46"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.64"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("E51").Value = "  +2.77%  "
